$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.072.49'
$ws.Range('E2').Value = '  -0.08%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.652.67'
$ws.Range('E3').Value = '  -0.46%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.27%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.70'
$ws.Range('E5').Value = '  +3.41%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5239'
$ws.Range('E6').Value = '  +1.76%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.000'
$ws.Range('E7').Value = '  -0.27%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2624'
$ws.Range('E8').Value = '  +1.76%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06372'
$ws.Range('E9').Value = '  +1.61%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.78'
$ws.Range('E10').Value = '  -0.63%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07737'
$ws.Range('E11').Value = '  +2.84%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.451'
$ws.Range('E12').Value = '  +1.21%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.639.39'
$ws.Range('E13').Value = '  -1.03%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.873.87'
$ws.Range('E14').Value = '  -0.56%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5511'
$ws.Range('E15').Value = '  +2.41%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₅8278'
$ws.Range('E16').Value = '  +4.99%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.95'
$ws.Range('E17').Value = '  -1.73%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '26.075.81'
$ws.Range('E18').Value = '  -0.11%  '
$ws.Range('E19').Value = '  -0.21%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.748'
$ws.Range('E20').Value = '  +1.35%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '190.41'
$ws.Range('E21').Value = '  +2.01%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.26'
$ws.Range('E22').Value = '  +1.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.351'
$ws.Range('E23').Value = '  +3.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.001'
$ws.Range('E24').Value = '  -0.22%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.21'
$ws.Range('E25').Value = '  -3.29%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1248'
$ws.Range('E26').Value = '  +3.32%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.402'
$ws.Range('E27').Value = '  +0.37%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.02'
$ws.Range('E28').Value = '  +2.83%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.418'
$ws.Range('E29').Value = '  +3.12%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05945'
$ws.Range('E30').Value = '  -2.62%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.262'
$ws.Range('E31').Value = '  +0.18%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.440'
$ws.Range('E32').Value = '  -0.72%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.417'
$ws.Range('E33').Value = '  +0.55%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.652'
$ws.Range('E34').Value = '  +1.68%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9966'
$ws.Range('E35').Value = '  +1.29%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.399'
$ws.Range('E36').Value = '  +0.49%  '
$ws.Range('E37').Value = '  +0.20%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5639'
$ws.Range('E38').Value = '  -3.79%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01603'
$ws.Range('E39').Value = '  +0.82%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.875'
$ws.Range('E40').Value = '  -1.56%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8568'
$ws.Range('E41').Value = '  +1.27%  '
$ws.Range('E42').Value = '  -0.18%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.028.01'
$ws.Range('E43').Value = '  -7.00%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '99.31'
$ws.Range('E44').Value = '  -0.75%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.798.23'
$ws.Range('E45').Value = '  -0.73%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0₈109'
$ws.Range('E46').Value = '  -1.74%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '55.84'
$ws.Range('E47').Value = '  +1.92%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.003'
$ws.Range('E48').Value = '  +0.19%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.046'
$ws.Range('E49').Value = '  +1.03%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05149'
$ws.Range('E50').Value = '  -1.60%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4210'
$ws.Range('E51').Value = '  -0.78%  '
